$wb = $excel.ActiveWorkbook

# Rename "InvalidNegotiationsFee" -> "InvalidTradingFees"
$wb.Worksheets.Item("InvalidNegotiationsFee").Name = "InvalidTradingFees"

# Rename "InvalidNegotiationFeesSummary" -> "InvalidTradingFeesSummary"
$wb.Worksheets.Item("InvalidNegotiationFeesSummary").Name = "InvalidTradingFeesSummary"

# Make the renamed "InvalidTradingFees" sheet the active sheet/tab
$wb.Worksheets.Item("InvalidTradingFees").Activate()
